$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) cells being updated. Some of the new values look like
# plain numbers (e.g. "232.28"), which Excel would otherwise auto-convert
# to a Number when assigned via .Value - but the workbook stores every
# Price cell as text (inline string), so we force a Text number format on
# each cell first (one at a time - applying it to a multi-area Range only
# affects the first area) and then put the style back to Normal afterwards
# so we do not leave a stray per-cell style behind.
$priceRefs = @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D21", "D22", "D23", "D25", "D26", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($ref in $priceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.763.53"
$ws.Range("D3").Value = "1.806.18"
$ws.Range("D5").Value = "232.28"
$ws.Range("D6").Value = "0.5915"
$ws.Range("D8").Value = "0.2780"
$ws.Range("D9").Value = "0.06833"
$ws.Range("D10").Value = "23.32"
$ws.Range("D11").Value = "0.07505"
$ws.Range("D12").Value = "1.804.09"
$ws.Range("D13").Value = "4.772"
$ws.Range("D14").Value = "0.6238"
$ws.Range("D15").Value = "2.050.74"
$ws.Range("D16").Value = "0.000009273"
$ws.Range("D17").Value = "75.82"
$ws.Range("D18").Value = "28.719.56"
$ws.Range("D21").Value = "211.65"
$ws.Range("D22").Value = "11.49"
$ws.Range("D23").Value = "6.843"
$ws.Range("D25").Value = "154.29"
$ws.Range("D26").Value = "7.885"
$ws.Range("D28").Value = "16.45"
$ws.Range("D29").Value = "1.426"
$ws.Range("D30").Value = "0.06173"
$ws.Range("D32").Value = "3.787"
$ws.Range("D33").Value = "3.765"
$ws.Range("D34").Value = "1.733"
$ws.Range("D35").Value = "1.066"
$ws.Range("D36").Value = "0.6436"
$ws.Range("D37").Value = "2.495"
$ws.Range("D38").Value = "2.720"
$ws.Range("D39").Value = "6.598"
$ws.Range("D40").Value = "0.01709"
$ws.Range("D41").Value = "1.143.50"
$ws.Range("D42").Value = "0.8828"
$ws.Range("D43").Value = "1.008"
$ws.Range("D44").Value = "100.30"
$ws.Range("D45").Value = "1.961.37"
$ws.Range("D46").Value = "60.48"
$ws.Range("D48").Value = "1.606"
$ws.Range("D49").Value = "8.361"
$ws.Range("D50").Value = "0.05473"
$ws.Range("D51").Value = "0.4487"

foreach ($ref in $priceRefs) {
    $ws.Range($ref).Style = "Normal"
}

# Volume(1h) (column E) cells being updated. These already round-trip as
# text because of the padding spaces around the percentage, so no special
# handling is required.
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("E9").Value = "  -3.23%  "
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("E16").Value = "  -6.70%  "
$ws.Range("E17").Value = "  -3.33%  "
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("E19").Value = "  -6.22%  "
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("E21").Value = "  -6.35%  "
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("E29").Value = "  -4.25%  "
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("E35").Value = "  -4.93%  "
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("E41").Value = "  -5.76%  "
$ws.Range("E42").Value = "  -2.57%  "
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("E46").Value = "  -3.24%  "
$ws.Range("E47").Value = "  -4.24%  "
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("E51").Value = "  -1.45%  "
